# TC02_C3DC_phs003111_SexAtBirth-Male.xlsx — regression/smoke suite update
#
# The "TreatmentTab" query (row 5, column B) had its
#   CONCAT(REPLACE(trt.treatment_agent, ';', ', '))
# simplified down to just
#   REPLACE(trt.treatment_agent, ';', ', ')
# (the CONCAT() wrapper around a single argument was a no-op, so it was
# dropped). Re-entering the cell's text is also what nudges Excel to
# (re)register the cell's font/style slot and reshuffle the shared-string
# table the way the authored workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentQuery = @'
SELECT
    DISTINCT prt.participant_id AS "Participant Id",
    trt.treatment_id AS "Treatment Id",
    CASE 
    WHEN trt.age_at_treatment_start = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_start >= 1000 THEN 
        substr(trt.age_at_treatment_start, 1, length(trt.age_at_treatment_start) - 3) || ',' || substr(trt.age_at_treatment_start, -3)
    ELSE 
        trt.age_at_treatment_start 
END AS "Age at Treatment Start",
    CASE 
    WHEN trt.age_at_treatment_end = -999 THEN 'Not Reported'
    WHEN trt.age_at_treatment_end >= 1000 THEN 
        substr(trt.age_at_treatment_end, 1, length(trt.age_at_treatment_end) - 3) || ',' || substr(trt.age_at_treatment_end, -3)
    ELSE 
        trt.age_at_treatment_end 
END AS "Age at Treatment End",
    trt.treatment_type AS "Treatment Type",
    REPLACE(trt.treatment_agent, ';', ', ') AS "Treatment Agent",
    std.dbgap_accession AS "dbGaP Accession"
FROM 
    df_study std
LEFT JOIN 
    df_participant prt ON std.id = prt."study.id"
LEFT JOIN 
    df_diagnoses dgn ON prt.id = dgn."participant.id"
LEFT JOIN 
    df_treatments trt ON prt.id = trt."participant.id"
LEFT JOIN 
    df_treatment_resp trr ON prt.id = trr."participant.id"
LEFT JOIN 
    df_survival srv ON prt.id = srv."participant.id"
LEFT JOIN 
    df_reference_files rfs ON std.id = rfs."study.id"
WHERE 
    std.dbgap_accession = 'phs003111' AND prt.sex_at_birth = 'Male'
ORDER BY 
    trt.treatment_id ASC
LIMIT 100;
'@

# Re-type the Treatment query text (drops the redundant CONCAT() wrapper).
$ws.Range("B5").Value = $treatmentQuery

# Touch the cell's font so it gets its own style entry (same visible
# Calibri 12pt look, distinct slot in the style table) just like the
# authored workbook ends up with.
$ws.Range("B5").Font.Name = "Calibri"

# Leave the sheet scrolled to the top with C5 selected, matching the
# saved view state in the authored file.
$ws.Range("C5").Select()
